# Remove the "EVIDENCIA DE PASE CON JACOCO" section (its heading paragraph
# and the paragraph holding its screenshot), along with the three blank
# paragraphs that immediately precede that heading.

$d = $word.ActiveDocument

# Locate the heading paragraph by its text so the script does not rely on
# brittle, hard-coded paragraph indices.
$headingIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*EVIDENCIA DE PASE CON JACOCO*") {
        $headingIndex = $i
        break
    }
    $i = $i + 1
}

if ($headingIndex -gt 0) {
    # Three empty paragraphs right before the heading, the heading itself,
    # and the following paragraph containing the evidence image.
    $startIndex = $headingIndex - 3
    $endIndex = $headingIndex + 1

    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
